# Refresh the cryptos list (prices + 1h volume %) per the upstream
# GitHub Actions data pull, including the VeChain / TrustWalletToken
# row swap at rows 39-40.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.774.86"
$ws.Range("E2").Value = "  -1.83%  "

$ws.Range("D3").Value = "'1.870.01"

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'300.15"
$ws.Range("E5").Value = "  -2.38%  "

$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").Value = "'0.5324"

$ws.Range("D8").Value = "'0.3729"
$ws.Range("E8").Value = "  -2.20%  "

$ws.Range("D9").Value = "'0.07156"
$ws.Range("E9").Value = "  -1.88%  "

$ws.Range("D10").Value = "'21.47"
$ws.Range("E10").Value = "  -2.38%  "

$ws.Range("D11").Value = "'0.8862"
$ws.Range("E11").Value = "  -1.81%  "

$ws.Range("D12").Value = "'0.08157"
$ws.Range("E12").Value = "  -0.06%  "

$ws.Range("D13").Value = "'1.869.31"
$ws.Range("E13").Value = "  +28.72%  "

$ws.Range("D14").Value = "'92.46"
$ws.Range("E14").Value = "  -3.71%  "

$ws.Range("D15").Value = "'5.289"
$ws.Range("E15").Value = "  -1.35%  "

$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = "  -0.04%  "

$ws.Range("D17").Value = "'14.83"
$ws.Range("E17").Value = "  +0.46%  "

$ws.Range("D18").Value = "'0.000008486"
$ws.Range("E18").Value = "  -2.12%  "

$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  -0.01%  "

$ws.Range("D20").Value = "'26.801.11"
$ws.Range("E20").Value = "  -1.84%  "

$ws.Range("D21").Value = "'4.970"
$ws.Range("E21").Value = "  -2.88%  "

$ws.Range("E22").Value = "  -1.90%  "

$ws.Range("D23").Value = "'6.377"
$ws.Range("E23").Value = "  -2.17%  "

$ws.Range("D24").Value = "'2.284"
$ws.Range("E24").Value = "  -0.97%  "

$ws.Range("D25").Value = "'145.87"
$ws.Range("E25").Value = "  -2.74%  "

$ws.Range("D26").Value = "'1.732"
$ws.Range("E26").Value = "  -0.10%  "

$ws.Range("D27").Value = "'18.02"
$ws.Range("E27").Value = "  -1.29%  "

$ws.Range("D28").Value = "'113.67"
$ws.Range("E28").Value = "  -2.79%  "

$ws.Range("D29").Value = "'4.698"
$ws.Range("E29").Value = "  -3.16%  "

$ws.Range("D30").Value = "'4.633"
$ws.Range("E30").Value = "  -4.40%  "

$ws.Range("D31").Value = "'0.09103"
$ws.Range("E31").Value = "  -1.67%  "

$ws.Range("D32").Value = "'0.8110"
$ws.Range("E32").Value = "  -2.26%  "

$ws.Range("D33").Value = "'0.05016"
$ws.Range("E33").Value = "  -1.05%  "

$ws.Range("D34").Value = "'1.174"
$ws.Range("E34").Value = "  -4.60%  "

$ws.Range("D35").Value = "'2.947"
$ws.Range("E35").Value = "  -1.80%  "

$ws.Range("D36").Value = "'0.6129"
$ws.Range("E36").Value = "  +5.65%  "

$ws.Range("D37").Value = "'2.651"
$ws.Range("E37").Value = "  -2.61%  "

$ws.Range("D38").Value = "'3.177"
$ws.Range("E38").Value = "  -5.17%  "

$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'1.070"
$ws.Range("E39").Value = "  -0.54%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.01944"
$ws.Range("E40").Value = "  -2.96%  "

$ws.Range("D41").Value = "'0.5284"
$ws.Range("E41").Value = "  +7.25%  "

$ws.Range("D42").Value = "'6.472"
$ws.Range("E42").Value = "  -1.78%  "

$ws.Range("D43").Value = "'8.724"
$ws.Range("E43").Value = "  -5.45%  "

$ws.Range("D44").Value = "'115.32"
$ws.Range("E44").Value = "  -1.25%  "

$ws.Range("D45").Value = "'0.1492"
$ws.Range("E45").Value = "  -2.15%  "

$ws.Range("D46").Value = "'1.001"
$ws.Range("E46").Value = "  +0.03%  "

$ws.Range("D47").Value = "'1.641"
$ws.Range("E47").Value = "  -0.44%  "

$ws.Range("D48").Value = "'9.941"
$ws.Range("E48").Value = "  -2.60%  "

$ws.Range("D49").Value = "'37.31"
$ws.Range("E49").Value = "  -4.39%  "

$ws.Range("D50").Value = "'0.06059"
$ws.Range("E50").Value = "  -1.75%  "

$ws.Range("D51").Value = "'62.08"
$ws.Range("E51").Value = "  -3.88%  "
